# MODEL-INPUT CHANGES for vl, saved
#
# The "Parameters" sheet had a row for the parameter "requiredvl"
# ("Number of VL tests recommended per person per year"). That row is
# removed entirely (not just cleared) so every row below it shifts up
# by one, and the two shared strings that were only used by that row
# (the label and the short code) drop out of the workbook once unused.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure we're editing the "Parameters" sheet regardless of which
# sheet happened to be active when the workbook was opened.
if ($ws.Name -ne "Parameters") {
    $ws = $wb.Worksheets.Item("Parameters")
    $ws.Activate()
}

# Row 70 holds "requiredvl" / "Number of VL tests recommended per
# person per year". Delete the entire row and shift the rows below it
# up, just like selecting row 70 and choosing Delete in Excel.
$ws.Rows.Item(70).Delete()

# Leave the selection on the row that is now in row 70's place (what
# used to be row 71), matching the post-delete selection state.
$ws.Range("A70:XFD70").Select()
